$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E2 value
$ws.Range("E2").Value = 11111003035

# Update selection to G2
$ws.Range("G2").Select()
